$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B78").Value = 5498504
$ws.Range("B79").Value = 5499423
$ws.Range("F78").Value = "NK Maribor"
$ws.Range("F79").Value = "Olimpija Ljubljana"
$ws.Range("G78").Value = "NK Bravo"
$ws.Range("G79").Value = "NK Celje"
$ws.Range("H78").Value = 1
$ws.Range("H79").Value = 0
$ws.Range("I78").Value = 1
$ws.Range("I79").Value = 2
$ws.Range("J78").Value = "D"
$ws.Range("J79").Value = "A"
$ws.Range("K78").Value = 1.571
$ws.Range("K79").Value = 2.5
$ws.Range("L78").Value = 3.8
$ws.Range("L79").Value = 3.3
$ws.Range("M78").Value = 4.75
$ws.Range("M79").Value = 2.5
$ws.Range("N78").Value = 1.533
$ws.Range("N79").Value = 2.55
$ws.Range("O78").Value = 4
$ws.Range("O79").Value = 3.25
$ws.Range("P78").Value = 4.75
$ws.Range("P79").Value = 2.45
$ws.Range("Q78").Value = -1
$ws.Range("Q79").Value = 0
$ws.Range("T78").Value = 2.75
$ws.Range("T79").Value = 2.5
$ws.Range("U78").Value = 1.875
$ws.Range("U79").Value = 1.85
$ws.Range("V78").Value = 1.925
$ws.Range("V79").Value = 1.95
$ws.Range("X78").Value = 3
$ws.Range("X79").Value = -1
$ws.Range("Y78").Value = -1
$ws.Range("Y79").Value = 1.45
$ws.Range("AC78").Value = 0.925
$ws.Range("AC79").Value = 0.95
$ws.Range("B82").Value = 6814327
$ws.Range("B83").Value = 6816473
$ws.Range("F82").Value = "NS Mura"
$ws.Range("F83").Value = "NK Bravo"
$ws.Range("G82").Value = "NK Domzale"
$ws.Range("G83").Value = "NK Rogaska"
$ws.Range("I82").Value = 3
$ws.Range("I83").Value = 0
$ws.Range("J82").Value = "A"
$ws.Range("J83").Value = "H"
$ws.Range("K82").Value = 2
$ws.Range("K83").Value = 1.8
$ws.Range("L82").Value = 3.3
$ws.Range("L83").Value = 3.5
$ws.Range("M82").Value = 3.4
$ws.Range("M83").Value = 4
$ws.Range("N82").Value = 1.909
$ws.Range("N83").Value = 2.05
$ws.Range("O82").Value = 3.4
$ws.Range("O83").Value = 3
$ws.Range("Q82").Value = -0.5
$ws.Range("Q83").Value = -0.25
$ws.Range("R82").Value = 1.95
$ws.Range("R83").Value = 1.75
$ws.Range("S82").Value = 1.85
$ws.Range("S83").Value = 2.05
$ws.Range("T82").Value = 2.5
$ws.Range("T83").Value = 2.25
$ws.Range("U82").Value = 1.9
$ws.Range("U83").Value = 1.95
$ws.Range("V82").Value = 1.9
$ws.Range("V83").Value = 1.85
$ws.Range("W82").Value = -1
$ws.Range("W83").Value = 1.05
$ws.Range("Y82").Value = 2.75
$ws.Range("Y83").Value = -1
$ws.Range("Z82").Value = -1
$ws.Range("Z83").Value = 0.75
$ws.Range("AA82").Value = 0.8500000000000001
$ws.Range("AA83").Value = -1
$ws.Range("AB82").Value = 0.8999999999999999
$ws.Range("AB83").Value = -0.5
$ws.Range("AC82").Value = -1
$ws.Range("AC83").Value = 0.425
$ws.Range("B89").Value = 6814328
$ws.Range("B90").Value = 6814330
$ws.Range("F89").Value = "NK Domzale"
$ws.Range("F90").Value = "NK Maribor"
$ws.Range("G89").Value = "NK Bravo"
$ws.Range("G90").Value = "NK Aluminij"
$ws.Range("I89").Value = 1
$ws.Range("I90").Value = 0
$ws.Range("J89").Value = "D"
$ws.Range("J90").Value = "H"
$ws.Range("K89").Value = 2.35
$ws.Range("K90").Value = 1.363
$ws.Range("L89").Value = 3.1
$ws.Range("L90").Value = 4.5
$ws.Range("M89").Value = 2.9
$ws.Range("M90").Value = 7
$ws.Range("N89").Value = 2.15
$ws.Range("N90").Value = 1.4
$ws.Range("O89").Value = 3.1
$ws.Range("O90").Value = 4.5
$ws.Range("P89").Value = 3.3
$ws.Range("P90").Value = 7
$ws.Range("Q89").Value = -0.25
$ws.Range("Q90").Value = -1.25
$ws.Range("R89").Value = 1.925
$ws.Range("R90").Value = 1.85
$ws.Range("S89").Value = 1.875
$ws.Range("S90").Value = 1.95
$ws.Range("T89").Value = 2.25
$ws.Range("T90").Value = 2.75
$ws.Range("U89").Value = 1.95
$ws.Range("U90").Value = 1.8
$ws.Range("V89").Value = 1.85
$ws.Range("V90").Value = 2
$ws.Range("W89").Value = -1
$ws.Range("W90").Value = 0.3999999999999999
$ws.Range("X89").Value = 2.1
$ws.Range("X90").Value = -1
$ws.Range("AA89").Value = 0.4375
$ws.Range("AA90").Value = 0.475
$ws.Range("AB89").Value = -0.5
$ws.Range("AB90").Value = -1
$ws.Range("AC89").Value = 0.425
$ws.Range("AC90").Value = 1
$ws.Range("N188").Value = 6.5
$ws.Range("O188").Value = 4.75
$ws.Range("P188").Value = 1.4
$ws.Range("R188").Value = 1.95
$ws.Range("S188").Value = 1.85
$ws.Range("U188").Value = 1.775
$ws.Range("V188").Value = 2.025
$ws.Range("U189").Value = 1.825
$ws.Range("V189").Value = 1.975
$ws.Range("R190").Value = 1.95
$ws.Range("S190").Value = 1.85
$ws.Range("U190").Value = 1.8
$ws.Range("V190").Value = 2
$ws.Range("U191").Value = 1.975
$ws.Range("V191").Value = 1.825
$ws.Range("N192").Value = 2.45
$ws.Range("P192").Value = 2.55
$ws.Range("R192").Value = 1.85
$ws.Range("S192").Value = 1.95
